$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data-wide-value")

# Insert a new column before column B, shifting existing B:F -> C:G
$ws.Columns.Item(2).Insert()

# Determine the last used row (data goes from row 1 header to row 112)
$lastRow = $ws.Cells(1,1).End(-4121).Row  # xlDown = -4121

# Header for the new column
$ws.Range("B1").Value = "budget-type"

# Fill the new column with the constant "budget" for every data row
$ws.Range("B2:B" + $lastRow).Value = "budget"
